$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set the new date header in AL1, matching the style (bold, General format) of the other header cells
$ws.Range("AL1").NumberFormat = "@"
$ws.Range("AL1").Value = "2020/05/27"
$ws.Range("AL1").NumberFormat = "General"
$ws.Range("AL1").Font.Bold = $true

# New column AL holds confirmed-case counts for 2020/05/27, one per concelho (row)
$newCounts = @{
    2 = 16
    3 = 65
    4 = 88
    5 = 76
    6 = 7
    7 = 8
    8 = 39
    9 = 23
    10 = 78
    11 = 6
    12 = 4
    13 = 361
    14 = 4
    15 = 20
    16 = 8
    17 = 12
    18 = 24
    19 = 728
    20 = 95
    21 = 78
    22 = 43
    23 = 5
    24 = 73
    25 = 8
    26 = 41
    27 = 12
    28 = 334
    29 = 71
    30 = 31
    31 = 299
    32 = 217
    33 = 5
    34 = 15
    35 = 35
    36 = 4
    37 = 1213
    38 = 119
    39 = 18
    40 = 9
    41 = 26
    42 = 4
    44 = 39
    45 = 19
    46 = 69
    47 = 6
    48 = 11
    49 = 57
    50 = 534
    51 = 7
    52 = 20
    53 = 104
    54 = 3
    55 = 6
    56 = 23
    57 = 9
    58 = 26
    59 = 24
    60 = 577
    61 = 155
    62 = 47
    63 = 7
    64 = 4
    65 = 9
    66 = 8
    67 = 89
    68 = 46
    69 = 95
    70 = 26
    71 = 121
    72 = 66
    73 = 403
    74 = 34
    76 = 4
    77 = 3
    78 = 26
    79 = 4
    80 = 10
    81 = 1079
    82 = 22
    83 = 11
    84 = 25
    85 = 710
    86 = 6
    87 = 131
    88 = 9
    89 = 4
    90 = 39
    91 = 83
    92 = 2254
    93 = 64
    94 = 894
    95 = 7
    96 = 20
    97 = 332
    98 = 22
    99 = 5
    100 = 115
    101 = 944
    102 = 78
    103 = 7
    104 = 90
    105 = 18
    106 = 1275
    107 = 18
    108 = 63
    109 = 6
    110 = 17
    111 = 10
    112 = 21
    113 = 5
    114 = 15
    115 = 146
    116 = 116
    119 = 7
    120 = 30
    121 = 116
    122 = 11
    123 = 71
    124 = 20
    125 = 17
    126 = 3
    127 = 19
    128 = 29
    130 = 9
    131 = 477
    132 = 390
    133 = 13
    134 = 210
    135 = 12
    136 = 39
    137 = 26
    138 = 28
    139 = 652
    140 = 308
    141 = 27
    142 = 334
    143 = 7
    144 = 3
    145 = 21
    146 = 170
    147 = 5
    148 = 10
    149 = 65
    150 = 23
    151 = 60
    152 = 8
    153 = 6
    154 = 6
    155 = 28
    157 = 6
    158 = 6
    159 = 38
    160 = 1349
    161 = 9
    162 = 4
    163 = 61
    164 = 157
    165 = 7
    166 = 64
    167 = 3
    168 = 19
    169 = 7
    170 = 11
    171 = 9
    172 = 5
    173 = 480
    174 = 3
    175 = 111
    176 = 17
    177 = 391
    178 = 3
    179 = 75
    180 = 9
    181 = 3
    182 = 8
    183 = 13
    184 = 335
    185 = 6
    186 = 14
    187 = 4
    188 = 33
    189 = 97
    190 = 43
    191 = 24
    192 = 1140
    193 = 4
    194 = 25
    195 = 35
    196 = 30
    197 = 10
    198 = 12
    199 = 13
    200 = 23
    201 = 18
    202 = 42
    203 = 18
    204 = 146
    205 = 22
    206 = 155
    207 = 12
    208 = 757
    209 = 6
    210 = 5
    211 = 8
    212 = 186
    213 = 37
    214 = 7
    215 = 292
    216 = 4
    217 = 377
    219 = 8
    220 = 399
    221 = 72
    222 = 1553
    223 = 6
    224 = 6
    225 = 153
    226 = 13
    227 = 233
    228 = 11
    229 = 30
    230 = 101
    231 = 137
    232 = 8
}

foreach ($row in $newCounts.Keys) {
    $ws.Cells.Item($row, 38).Value = $newCounts[$row]
}
